# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 86d68523-e520-460c-9dab-5e3c018e5b46 row on both the
# zh-cn and de-de report sheets (row 4, columns E and H).

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-24 12:58:40"
$wsZh.Range("H4").Value = "2016-03-24 12:59:28"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-24 12:58:46"
$wsDe.Range("H4").Value = "2016-03-24 12:59:35"
